# Fix student processing bug: Student S0002 (Jos van Weert, row 3 on the
# "Student" sheet) was linked to assignment "O0003", which does not exist
# in the Project sheet (S0002 -> P0002). Correct it to "O0002".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student")
$ws.Range("B3").Value = "O0002"
